$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / changed date) for rows 2-9 from 45221 to 45224
$ws.Range("C2:C9").Value = 45224
